$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A243").Value = 'What''s the maximum curves I can load in a data file?'
$ws.Range("B243").Value = 'llama3.2:latest'
$ws.Range("C243").Value = 'According to Document 30, the maximum number of data files you can load is unlimited. However, it does not specify the maximum number of curves per data file.
To answer your question accurately, let''s look at another document that might provide more information on this topic.
Document 29 states: "Deviation from the above may cause errors during loading data from files." It also mentions that if possible, you should ask for LAS files from your data provider. However, it does not specify a maximum number of curves per file.
Another relevant document is Document 27, which states: "Curve shade name length" with a value of 20. This might imply that the maximum curve name length is 20 characters, but it doesn''t provide information on the maximum number of curves you can load in a single data file.
Based on these documents, I couldn''t find any specific information on the maximum number of curves you can load in a data file.'

$ws.Range("A244").Value = 'What''s the maximum curves I can load in a data file?'
$ws.Range("B244").Value = 'llama3.2:latest'
$ws.Range("C244").Value = 'Based on Document 29, which states: "Deviation from the above may cause errors during loading data from files." and considering that the curves are being loaded into columns in Geo so the maximum number of curves correspond to the maximum number of columns in the data file, the answer is:
You can load up to 200 tracks (curves) per data file.'

$ws.Range("A245").Value = 'What''s the maximum curves I can load in a data file?'
$ws.Range("B245").Value = 'llama3.2:latest'
$ws.Range("C245").Value = 'Based on Document 29, which states: "Deviation from the above may cause errors during loading data from files." and considering that the curves are being loaded into columns in Geo so the maximum number of curves correspond to the maximum number of columns in the data file, the answer is:
You can load up to 200 tracks (curves) per data file.'

$ws.Range("A246").Value = 'What''s the maximum curves I can load in a data file?'
$ws.Range("B246").Value = 'llama3.2:latest'
$ws.Range("C246").Value = 'Based on Document 29, which states: "Deviation from the above may cause errors during loading data from files." and considering that the curves are being loaded into columns in Geo so the maximum number of curves correspond to the maximum number of columns in the data file, the answer is:
You can load up to 200 tracks (curves) per data file.'

$ws.Range("A247").Value = 'What is the limit to the number of columns per data file?'
$ws.Range("B247").Value = 'llama3.2:latest'
$ws.Range("C247").Value = 'Based on Document 31, which states: "Columns per data file", the maximum number of columns that can be loaded in a GEO file value is 450.'

$ws.Range("A248").Value = 'Can curves be loaded into columns?'
$ws.Range("B248").Value = 'llama3.2:latest'
$ws.Range("C248").Value = 'Yes, curves can be loaded into columns in GEO. To do this, select the table containing the data you want to display as a curve, then click "Select Curves..." and enter a Curve Mnemonic for each Table Column you wish to generate a curve from. The maximum number of columns that can be loaded is 450.'

$ws.Range("A249").Value = 'What''s the maximum curves I can load in a data file?'
$ws.Range("B249").Value = 'llama3.2:latest'
$ws.Range("C249").Value = 'Based on Document 29, which states: "Deviation from the above may cause errors during loading data from files." and considering that the curves are being loaded into columns in Geo so the maximum number of curves correspond to the maximum number of columns in the data file, the answer is:
You can load up to 200 tracks (curves) per data file.'

$ws.Range("A250").Value = 'What''s the maximum number of curves I can load in a data file?'
$ws.Range("B250").Value = 'llama3.2:latest'
$ws.Range("C250").Value = 'Based on Document 29, which states: "Deviation from the above may cause errors during loading data from files." and considering that the curves are being loaded into columns in Geo so the maximum number of curves correspond to the maximum number of columns in the data file, the answer is:
You can load up to 200 tracks (curves) per data file.'

$ws.Range("A251").Value = 'What''s the maximum number of curves I can load?'
$ws.Range("B251").Value = 'llama3.2:latest'
$ws.Range("C251").Value = 'Based on Document 29, which states: "Deviation from the above may cause errors during loading data from files." and considering that the curves are being loaded into columns in Geo so the maximum number of curves correspond to the maximum number of columns in the data file, the answer is: You can load up to 200 tracks (curves) per data file.'

$ws.Range("A252").Value = 'What''s the maximum number of curves I can load?'
$ws.Range("B252").Value = 'llama3.2:latest'
$ws.Range("C252").Value = 'Based on Document 29, which states: "Deviation from the above may cause errors during loading data from files." and considering that the curves are being loaded into columns in Geo so the maximum number of curves correspond to the maximum number of columns in the data file, the answer is: You can load up to 450 tracks (curves) per data file.'

$ws.Range("A253").Value = 'Question 10: What''s the maximum number of curves I can load?'
$ws.Range("B253").Value = 'llama3.2:latest'
$ws.Range("C253").Value = 'Based on Document 29, which states: "Deviation from the above may cause errors during loading data from files." and considering that the curves are being loaded into columns in Geo so the maximum number of curves correspond to the maximum number of columns in the data file, the answer is: You can load up to 450 curves per data file.'

$ws.Range("A254").Value = 'What''s the maximum curves I can load in a data file?'
$ws.Range("B254").Value = 'llama3.2:latest'
$ws.Range("C254").Value = 'Based on Document 29, which states: "Deviation from the above may cause errors during loading data from files." and considering that the curves are being loaded into columns in Geo so the maximum number of curves correspond to the maximum number of columns in the data file, the answer is:
You can load up to 450 curves per data file.'
